$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.963.38'
$ws.Range("E2").Value = '  -0.23%  '
$ws.Range("D3").Value = '2.418.27'
$ws.Range("E3").Value = '  +0.09%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '''562.72'
$ws.Range("E5").Value = '  +0.93%  '
$ws.Range("D6").Value = '''142.78'
$ws.Range("E6").Value = '  -0.07%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '''0.530'
$ws.Range("E8").Value = '  -0.21%  '
$ws.Range("E9").Value = '  +0.14%  '
$ws.Range("E10").Value = '  -1.10%  '
$ws.Range("D11").Value = '''5.20'
$ws.Range("E11").Value = '  -3.90%  '
$ws.Range("D12").Value = '''0.350'
$ws.Range("E12").Value = '  -0.83%  '
$ws.Range("D13").Value = '''25.88'
$ws.Range("E13").Value = '  -0.86%  '
$ws.Range("D14").Value = '''0.0000173'
$ws.Range("E14").Value = '  -1.45%  '
$ws.Range("D15").Value = '2.854.33'
$ws.Range("E15").Value = '  -0.18%  '
$ws.Range("D16").Value = '61.806.55'
$ws.Range("E16").Value = '  -0.14%  '
$ws.Range("D17").Value = '2.421.12'
$ws.Range("E17").Value = '  +0.43%  '
$ws.Range("D18").Value = '''11.31'
$ws.Range("E18").Value = '  +1.55%  '
$ws.Range("D19").Value = '''323.83'
$ws.Range("E19").Value = '  -0.08%  '
$ws.Range("B20").Value = 'Polkadot'
$ws.Range("C20").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D20").Value = '''4.13'
$ws.Range("E20").Value = '  -1.51%  '
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").Value = '''6.82'
$ws.Range("E21").Value = '  +1.10%  '
$ws.Range("E22").Value = '  +0.01%  '
$ws.Range("D23").Value = '''66.63'
$ws.Range("E23").Value = '  +1.92%  '
$ws.Range("D24").Value = '''1.72'
$ws.Range("E24").Value = '  +0.37%  '
$ws.Range("D25").Value = '''8.68'
$ws.Range("E25").Value = '  -3.82%  '
$ws.Range("D26").Value = '''552.95'
$ws.Range("E26").Value = '  -5.71%  '
$ws.Range("D27").Value = '2.537.50'
$ws.Range("E27").Value = '  +0.64%  '
$ws.Range("E28").Value = '  -0.10%  '
$ws.Range("D29").Value = '0.0₃0930'
$ws.Range("E29").Value = '  -0.62%  '
$ws.Range("D30").Value = '''8.16'
$ws.Range("E30").Value = '  -1.41%  '
$ws.Range("D31").Value = '''1.39'
$ws.Range("E31").Value = '  -4.45%  '
$ws.Range("D32").Value = '''0.147'
$ws.Range("E32").Value = '  -1.18%  '
$ws.Range("E33").Value = '  -0.71%  '
$ws.Range("D34").Value = '''1.50'
$ws.Range("E34").Value = '  -4.16%  '
$ws.Range("E35").Value = '  -0.02%  '
$ws.Range("D36").Value = '''4.74'
$ws.Range("E36").Value = '  -0.75%  '
$ws.Range("E37").Value = '  -1.35%  '
$ws.Range("D38").Value = '''153.43'
$ws.Range("E38").Value = '  +2.21%  '
$ws.Range("D39").Value = '''5.43'
$ws.Range("E39").Value = '  -4.76%  '
$ws.Range("D40").Value = '''18.51'
$ws.Range("E40").Value = '  -0.92%  '
$ws.Range("E41").Value = '  -1.50%  '
$ws.Range("D42").Value = '''0.991'
$ws.Range("E42").Value = '  -0.82%  '
$ws.Range("D43").Value = '''146.84'
$ws.Range("E43").Value = '  -2.75%  '
$ws.Range("D44").Value = '''2.23'
$ws.Range("E44").Value = '  -4.66%  '
$ws.Range("D45").Value = '''3.63'
$ws.Range("E45").Value = '  -0.60%  '
$ws.Range("D46").Value = '''0.0526'
$ws.Range("E46").Value = '  -2.60%  '
$ws.Range("D47").Value = '''19.85'
$ws.Range("E47").Value = '  -1.83%  '
$ws.Range("D48").Value = '''0.591'
$ws.Range("E48").Value = '  +0.18%  '
$ws.Range("D49").Value = '''0.0919'
$ws.Range("E49").Value = '  -0.46%  '
$ws.Range("E50").Value = '  -1.05%  '
$ws.Range("E51").Value = '  +0.59%  '

# Reset style on cells that were entered with a leading apostrophe (forced-text marker)
# so no residual "quote prefix" cell formatting lingers on them.
foreach ($addr in @("D5","D6","D8","D11","D12","D13","D14","D18","D19","D20","D21","D23","D24","D25","D26","D30","D31","D32","D34","D36","D38","D39","D40","D42","D43","D44","D45","D46","D47","D48","D49")) {
    $ws.Range($addr).Style = "Normal"
}
